$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 32

# Columns A-D are textual (date/time/weekday/week) — force text type so
# Excel's COM layer doesn't auto-convert them to a date serial / number,
# matching the inlineStr cells used for every other row in this sheet.
function Set-TextCell($r, $c, $val) {
    $cell = $ws.Cells.Item($r, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextCell $row 1 "2025-01-18"
Set-TextCell $row 2 "12:10:45"
Set-TextCell $row 3 "Saturday"
Set-TextCell $row 4 "02"

# Columns E-T are plain numeric city figures.
$ws.Cells.Item($row, 5).Value  = 126841
$ws.Cells.Item($row, 6).Value  = 142140
$ws.Cells.Item($row, 7).Value  = 169233
$ws.Cells.Item($row, 8).Value  = 158537
$ws.Cells.Item($row, 9).Value  = -1
$ws.Cells.Item($row, 10).Value = 142937
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 192266
$ws.Cells.Item($row, 14).Value = 115576
$ws.Cells.Item($row, 15).Value = 45465
$ws.Cells.Item($row, 16).Value = 28503
$ws.Cells.Item($row, 17).Value = 65649
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 49066
$ws.Cells.Item($row, 20).Value = -1
